$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D store plain-looking numeric strings (e.g. "250.46") as
# TEXT in the source workbook (t="inlineStr"), matching the scraped price
# format (note some use "." as a thousands separator, e.g. "37.185.67").
# Force NumberFormat to Text first so Excel does not auto-coerce these into
# real numbers when we assign the new price strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.185.67"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.058.77"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "250.46"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +10.81%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.388"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").Value = "0.0795"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").Value = "16.34"
$ws.Range("E12").Value = "  +8.24%  "
$ws.Range("D13").Value = "2.356.53"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("E15").Value = "  +9.26%  "
$ws.Range("D16").Value = "2.056.64"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "18.25"
$ws.Range("E17").Value = "  +28.06%  "
$ws.Range("D18").Value = "37.138.81"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "75.67"
$ws.Range("E19").Value = "  +4.33%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -4.88%  "
$ws.Range("D21").Value = "5.46"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").Value = "239.43"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  +12.40%  "
$ws.Range("D26").Value = "169.67"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "9.46"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").Value = "20.12"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("E30").Value = "  +10.78%  "
$ws.Range("D31").Value = "4.84"
$ws.Range("E31").Value = "  +6.24%  "
$ws.Range("D32").Value = "0.0624"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +4.47%  "
$ws.Range("D34").Value = "0.0889"
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "2.28"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "1.73"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "0.110"
$ws.Range("E38").Value = "  +3.60%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "5.30"
$ws.Range("E40").Value = "  +31.41%  "
$ws.Range("D41").Value = "3.17"
$ws.Range("E41").Value = "  +14.39%  "
$ws.Range("D42").Value = "18.25"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "98.03"
$ws.Range("D46").Value = "2.51"
$ws.Range("E46").Value = "  +3.53%  "
$ws.Range("D47").Value = "1.298.10"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "2.246.38"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "3.53"
$ws.Range("E51").Value = "  -16.35%  "
